# Periode 4: fix a missing column in the Disposition formulas, then set the
# correct input value for Produktionsprogramm P3 (Periode 1 column).
$wb = $excel.ActiveWorkbook

# --- 1. Produktionsprogramm: P3 / Periode 1 changes from 200 to 50 ---
$wsProg = $wb.Worksheets.Item("Produktionsprogramm")
$wsProg.Range("C7").Value = 50

# --- 2. Fix the "missing column" bug: formulas in column O on the three
#        Disposition sheets didn't include column E (Rückstand/backlog
#        carried in via K of the previous block). Add "+E#" to row 5 and to
#        the shared formula anchored at row 6 (O6:O18) on each sheet. ---
$sheetNames = @("Disposition_P1", "Disposition_P2", "Disposition_P3")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("O5").Formula = "=IF((C5+E5+G5-I5-K5-M5)<0,0,C5+E5+G5-I5-K5-M5)"
    $ws.Range("O6:O18").Formula = "=IF((C6+E6+G6-I6-K6-M6)<0,0,C6+E6+G6-I6-K6-M6)"
}

$wb.Application.Calculate()

# --- 3. Restore selection / active sheet state to match the saved file ---
$wsProg.Range("C8").Select()

$wsP1 = $wb.Worksheets.Item("Disposition_P1")
$wsP1.Range("O5").Select()

$wsP2 = $wb.Worksheets.Item("Disposition_P2")
$wsP2.Range("O5:O18").Select()

$wsP3 = $wb.Worksheets.Item("Disposition_P3")
$wsP3.Range("R11").Select()
$wsP3.Activate()
